$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("earnings_debt")

# Row 2 (header row "Macedonia" / company_name "1")
$ws.Range("K2").Value = -4.69
$ws.Range("U2").Value = 2.93
$ws.Range("V2").Value = 0.09361022364217253
$ws.Range("W2").Value = -0.6020539152759949
$ws.Range("X2").Value = 0.1099973702347115
$ws.Range("Y2").Value = -0.7120512855107064
$ws.Range("AA2").Value = -0.2084084084084084
$ws.Range("AB2").Value = 0.116108343482846
$ws.Range("AC2").Value = -0.3245167518912544
$ws.Range("AD2").Value = 16.7
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 16.7
$ws.Range("AG2").Value = 13.77
$ws.Range("AH2").Value = 0.3479166666666667
$ws.Range("AI2").Value = 0.7844058243306716
$ws.Range("AJ2").Value = 0.3055247392944309
$ws.Range("AK2").Value = 0.75
$ws.Range("AL2").Value = 1.81
$ws.Range("AM2").Value = 1.748
$ws.Range("AN2").Value = -4.854651162790698
$ws.Range("AO2").Value = -1.917127071823205
$ws.Range("AP2").Value = -4.002906976744186
$ws.Range("AQ2").Value = -1.98512585812357

# Row 3 (Euromax Resources Ltd.)
$ws.Range("K3").Value = -4.69
$ws.Range("U3").Value = 2.93
$ws.Range("V3").Value = 0.09361022364217253
$ws.Range("W3").Value = -0.6020539152759949
$ws.Range("X3").Value = 0.1099973702347115
$ws.Range("Y3").Value = -0.7120512855107064
$ws.Range("AA3").Value = -0.2084084084084084
$ws.Range("AB3").Value = 0.116108343482846
$ws.Range("AC3").Value = -0.3245167518912544
$ws.Range("AD3").Value = 16.7
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 16.7
$ws.Range("AG3").Value = 13.77
$ws.Range("AH3").Value = 0.3479166666666667
$ws.Range("AI3").Value = 0.7844058243306716
$ws.Range("AJ3").Value = 0.3055247392944309
$ws.Range("AK3").Value = 0.75
$ws.Range("AL3").Value = 1.81
$ws.Range("AM3").Value = 1.748
$ws.Range("AN3").Value = -4.854651162790698
$ws.Range("AO3").Value = -1.917127071823205
$ws.Range("AP3").Value = -4.002906976744186
$ws.Range("AQ3").Value = -1.98512585812357
